$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "FilesTab" query (cell B4) still filtered on experimental_strategies:
# ["RNA-Seq"]. Reset it to an empty filter, matching the other tabs
# (ParticipantsTab / SamplesTab) which already use experimental_strategies: [].
$cell = $ws.Range("B4")
$query = $cell.Value2
$query = $query.Replace('experimental_strategies: ["RNA-Seq"],', 'experimental_strategies: [],')
$cell.Value2 = $query

# Row 4 auto-wraps this (still long) multi-line text, so keep it at Excel's
# maximum row height (the same effective height the row already had).
$ws.Rows.Item(4).RowHeight = 409.5

# Move the active selection from D4 to B4.
[void]$ws.Range("B4").Select()
